$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, per the diff
$changes = @{
    "E3"  = 12.79689999999999
    "A4"  = -21.37440000000002
    "C4"  = -11.17599999999999
    "E4"  = 12.4406
    "C5"  = -14.57100000000001
    "A6"  = -21.56830000000001
    "A7"  = -21.505
    "C8"  = -12.0579
    "E9"  = 13.10330000000001
    "E11" = 13.3214
    "E14" = 12.43270000000001
    "A16" = -20.48539999999998
    "C16" = -12.18420000000001
    "E18" = 13.1087
    "A20" = -22.94010000000002
    "C22" = -10.87619999999999
    "E25" = 13.25849999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
